$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Q-column anomaly flags to 0 ---
$ws.Range("Q61").Value = 0
$ws.Range("Q64").Value = 0
$ws.Range("Q70").Value = 0
$ws.Range("Q77").Value = 0

# --- Update O152 ---
$ws.Range("O152").Value = 2

# --- Append new weekly rows 156-179 ---
$newRows = @(
    @(45474, 199.7100067138672, 213.9499969482422, 194.5599975585938, 207.5700073242188, 207.5700073242188, 240399900, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(45481, 208.25, 223.25, 205.7400054931641, 222.4700012207031, 222.4700012207031, 178774698, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(45488, 225, 232, 210.25, 218.8899993896484, 218.8899993896484, 193679626, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(45495, 214, 226.6199951171875, 207, 224.6199951171875, 224.6199951171875, 164468500, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 226.6199951171875, 278.7000122070312, 224.1000061035156, 262.3399963378906, 262.3399963378906, 581982221, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(45509, 254, 275, 247.1999969482422, 267.0899963378906, 267.0899963378906, 537464038, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(45516, 263.5, 266.489990234375, 255.5, 264.4299926757812, 264.4299926757812, 144309902, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(45523, 278, 280.8999938964844, 256.010009765625, 262.6499938964844, 262.6499938964844, 427552275, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 263.7999877929688, 264.4200134277344, 248.5099945068359, 250.5299987792969, 250.5299987792969, 184529581, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(45537, 252, 262.1499938964844, 240.3999938964844, 259.9500122070312, 259.9500122070312, 372462299, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(45544, 260, 286.4500122070312, 256.2999877929688, 272.8999938964844, 272.8999938964844, 432919978, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0),
    @(45551, 274.5, 293, 269.2000122070312, 290.5, 290.5, 244174894, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(45558, 292, 298.25, 273.5, 278.1499938964844, 278.1499938964844, 235162814, 2024, 9, 23, 0, 0, 0, 39, 1, 0, 0),
    @(45565, 277.8999938964844, 277.8999938964844, 264.1000061035156, 275.2999877929688, 275.2999877929688, 172981785, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(45572, 278.5499877929688, 290, 261.6000061035156, 277.5499877929688, 277.5499877929688, 239339411, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(45579, 280.8999938964844, 285.25, 255.25, 257.4500122070312, 257.4500122070312, 223322675, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(45586, 258, 270.8999938964844, 242.1000061035156, 253.8000030517578, 253.8000030517578, 447063925, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0),
    @(45593, 252, 259.25, 240.3999938964844, 248.9900054931641, 248.9900054931641, 169414140, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0),
    @(45600, 242.9299926757812, 262.4500122070312, 239.4499969482422, 248.7299957275391, 248.7299957275391, 258639387, 2024, 11, 4, 0, 0, 0, 45, 2, 0, 2),
    @(45607, 244.75, 271, 244, 269.6600036621094, 269.6600036621094, 241394033, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0),
    @(45614, 274.2999877929688, 275.489990234375, 262.1099853515625, 264.2000122070312, 264.2000122070312, 187161112, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0),
    @(45621, 273.25, 289.7000122070312, 271.1700134277344, 279.760009765625, 279.760009765625, 357552407, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0),
    @(45628, 281.8999938964844, 304.6499938964844, 279, 302.9500122070312, 302.9500122070312, 262227530, 2024, 12, 2, 0, 0, 0, 49, 0, 0, 0),
    @(45635, 304.7000122070312, 304.7000122070312, 278.2000122070312, 288.25, 288.25, 207203572, 2024, 12, 9, 0, 0, 0, 50, 0, 0, 0)
)

$startRow = 156
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = [double]$row[0]
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat
    for ($j = 1; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = [double]$row[$j]
    }
    # Column R is left blank (empty string cell), matching the source data shape;
    # copy an existing blank "backup" cell to materialize an empty cell in column R.
    $ws.Range("R155").Copy($ws.Cells.Item($r, 18))
}

# --- Fill previously-blank backup cells (R154, R155) with numeric 0 ---
# (done after copying from R155 above so the copy source is still blank)
$ws.Range("R154").Value = 0
$ws.Range("R155").Value = 0
